$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 557.6923
$ws.Range("I6").Value = 465.55554
$ws.Range("J6").Value = 765
$ws.Range("K6").Value = 1396.66662
$ws.Range("L6").Value = 2295
$ws.Range("M6").Value = -1284.66662
$ws.Range("N6").Value = -2519
# Row 8
$ws.Range("H8").Value = 198.5
$ws.Range("I8").Value = 198.5
$ws.Range("K8").Value = 595.5
$ws.Range("M8").Value = -456.5
# Row 51
$ws.Range("H51").Value = 4816.6665
$ws.Range("J51").Value = 6125
$ws.Range("L51").Value = 6125
$ws.Range("N51").Value = -7093
# Row 64
$ws.Range("H64").Value = 43063.24
$ws.Range("I64").Value = 3128.8572
$ws.Range("K64").Value = 3128.8572
$ws.Range("M64").Value = -2880.8572
# Row 67
$ws.Range("H67").Value = 43063.24
$ws.Range("I67").Value = 3128.8572
$ws.Range("K67").Value = 3128.8572
$ws.Range("M67").Value = -2270.8572
# Row 76
$ws.Range("H76").Value = 3375.2424
$ws.Range("I76").Value = 3332.4644
$ws.Range("J76").Value = 3614.8
$ws.Range("K76").Value = 3332.4644
$ws.Range("L76").Value = 3614.8
$ws.Range("M76").Value = -3017.4644
$ws.Range("N76").Value = -4244.8
# Row 79
$ws.Range("H79").Value = 3375.2424
$ws.Range("I79").Value = 3332.4644
$ws.Range("J79").Value = 3614.8
$ws.Range("K79").Value = 3332.4644
$ws.Range("L79").Value = 3614.8
$ws.Range("M79").Value = -2240.4644
$ws.Range("N79").Value = -5798.8
# Row 129
$ws.Range("H129").Value = 631.25
$ws.Range("I129").Value = 276.9
$ws.Range("J129").Value = 828.1111
$ws.Range("K129").Value = 830.6999999999999
$ws.Range("L129").Value = 2484.3333
$ws.Range("M129").Value = 4169.3
$ws.Range("N129").Value = -12484.3333
# Row 137
$ws.Range("H137").Value = 3990.8445
$ws.Range("I137").Value = 4840.9414
$ws.Range("J137").Value = 3474.7144
$ws.Range("K137").Value = 14522.8242
$ws.Range("L137").Value = 10424.1432
$ws.Range("M137").Value = -11972.8242
$ws.Range("N137").Value = -15524.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3240.8
$ws.Range("I63").Value = 3118.3333
$ws.Range("K63").Value = 3118.3333
$ws.Range("M63").Value = -2432.3333
# Row 66
$ws.Range("H66").Value = 3240.8
$ws.Range("I66").Value = 3118.3333
$ws.Range("K66").Value = 15591.6665
$ws.Range("M66").Value = -12159.6665

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2006.6666
$ws.Range("I105").Value = 2006.6666
$ws.Range("K105").Value = 2006.6666
$ws.Range("M105").Value = -259.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 1593.6428
$ws.Range("I62").Value = 1467.2222
$ws.Range("J62").Value = 1821.2
$ws.Range("K62").Value = 1467.2222
$ws.Range("L62").Value = 1821.2
$ws.Range("M62").Value = -843.2221999999999
$ws.Range("N62").Value = -3069.2
# Row 65
$ws.Range("H65").Value = 1593.6428
$ws.Range("I65").Value = 1467.2222
$ws.Range("J65").Value = 1821.2
$ws.Range("K65").Value = 7336.111
$ws.Range("L65").Value = 9106
$ws.Range("M65").Value = -4216.111
$ws.Range("N65").Value = -15346
# Row 122
$ws.Range("H122").Value = 4115.727
$ws.Range("I122").Value = 3689.6086
$ws.Range("J122").Value = 5095.8
$ws.Range("K122").Value = 11068.8258
$ws.Range("L122").Value = 15287.4
$ws.Range("M122").Value = -8618.825800000001
$ws.Range("N122").Value = -20187.4

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 242.44444
$ws.Range("I4").Value = 320
$ws.Range("J4").Value = 180.4
$ws.Range("K4").Value = 960
$ws.Range("L4").Value = 541.2
$ws.Range("M4").Value = -848
$ws.Range("N4").Value = -765.2
# Row 5
$ws.Range("H5").Value = 20749.863
$ws.Range("I5").Value = 364.9643
$ws.Range("J5").Value = 45566.26
$ws.Range("K5").Value = 1094.8929
$ws.Range("L5").Value = 136698.78
$ws.Range("M5").Value = -982.8928999999998
$ws.Range("N5").Value = -136922.78
# Row 48
$ws.Range("H48").Value = 111113450
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 111113450
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -333340850
# Row 75
$ws.Range("H75").Value = 2139.4443
$ws.Range("I75").Value = 1197.5
$ws.Range("J75").Value = 2893
$ws.Range("K75").Value = 3592.5
$ws.Range("L75").Value = 8679
$ws.Range("M75").Value = -2594.5
$ws.Range("N75").Value = -10675
# Row 78
$ws.Range("H78").Value = 2139.4443
$ws.Range("I78").Value = 1197.5
$ws.Range("J78").Value = 2893
$ws.Range("K78").Value = 10777.5
$ws.Range("L78").Value = 26037
$ws.Range("M78").Value = -5785.5
$ws.Range("N78").Value = -36021
# Row 129
$ws.Range("H129").Value = 4445601.5
$ws.Range("I129").Value = 524.44446
$ws.Range("J129").Value = 6945957.5
$ws.Range("K129").Value = 1573.33338
$ws.Range("L129").Value = 20837872.5
$ws.Range("M129").Value = 3426.66662
$ws.Range("N129").Value = -20847872.5
# Row 135
$ws.Range("H135").Value = 20749.863
$ws.Range("I135").Value = 364.9643
$ws.Range("J135").Value = 45566.26
$ws.Range("K135").Value = 3284.6787
$ws.Range("L135").Value = 410096.34
$ws.Range("M135").Value = -749.6786999999999
$ws.Range("N135").Value = -415166.34

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4091.8518
$ws.Range("I70").Value = 4096
$ws.Range("K70").Value = 4096
$ws.Range("M70").Value = -3826
# Row 73
$ws.Range("H73").Value = 4091.8518
$ws.Range("I73").Value = 4096
$ws.Range("K73").Value = 4096
$ws.Range("M73").Value = -3160

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1963.75
$ws.Range("I7").Value = 1480
$ws.Range("J7").Value = 2770
$ws.Range("K7").Value = 1480
$ws.Range("L7").Value = 2770
$ws.Range("M7").Value = -1368
$ws.Range("N7").Value = -2994
# Row 16
$ws.Range("H16").Value = 1239.5264
$ws.Range("I16").Value = 797.1539
$ws.Range("J16").Value = 2198
$ws.Range("K16").Value = 797.1539
$ws.Range("L16").Value = 2198
$ws.Range("M16").Value = -627.1539
$ws.Range("N16").Value = -2538
# Row 22
$ws.Range("H22").Value = 1143.8
$ws.Range("I22").Value = 753.6
$ws.Range("J22").Value = 1534
$ws.Range("K22").Value = 753.6
$ws.Range("L22").Value = 1534
$ws.Range("M22").Value = -458.6
$ws.Range("N22").Value = -2124
# Row 27
$ws.Range("H27").Value = 1143.8
$ws.Range("I27").Value = 753.6
$ws.Range("J27").Value = 1534
$ws.Range("K27").Value = 753.6
$ws.Range("L27").Value = 1534
$ws.Range("M27").Value = -646.6
$ws.Range("N27").Value = -1748
# Row 40
$ws.Range("H40").Value = 2478.9678
$ws.Range("I40").Value = 2257.92
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 2257.92
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -2121.92
$ws.Range("N40").Value = -3672
# Row 68
$ws.Range("H68").Value = 2083.394
$ws.Range("I68").Value = 748.25
$ws.Range("J68").Value = 2663.8914
$ws.Range("K68").Value = 748.25
$ws.Range("L68").Value = 2663.8914
$ws.Range("M68").Value = 0.75
$ws.Range("N68").Value = -4161.8914
# Row 71
$ws.Range("H71").Value = 2083.394
$ws.Range("I71").Value = 748.25
$ws.Range("J71").Value = 2663.8914
$ws.Range("K71").Value = 3741.25
$ws.Range("L71").Value = 13319.457
$ws.Range("M71").Value = 2.75
$ws.Range("N71").Value = -20807.457
# Row 82
$ws.Range("H82").Value = 2164.5789
$ws.Range("I82").Value = 1356.6666
$ws.Range("J82").Value = 2691.4783
$ws.Range("K82").Value = 1356.6666
$ws.Range("L82").Value = 2691.4783
$ws.Range("M82").Value = -995.6666
$ws.Range("N82").Value = -3413.4783
# Row 85
$ws.Range("H85").Value = 2164.5789
$ws.Range("I85").Value = 1356.6666
$ws.Range("J85").Value = 2691.4783
$ws.Range("K85").Value = 1356.6666
$ws.Range("L85").Value = 2691.4783
$ws.Range("M85").Value = -108.6666
$ws.Range("N85").Value = -5187.478300000001
# Row 122
$ws.Range("H122").Value = 3059.5
$ws.Range("I122").Value = 3043
$ws.Range("J122").Value = 3101.6667
$ws.Range("K122").Value = 9129
$ws.Range("L122").Value = 9305.000100000001
$ws.Range("M122").Value = -6679
$ws.Range("N122").Value = -14205.0001
# Row 126
$ws.Range("H126").Value = 1963.75
$ws.Range("I126").Value = 1480
$ws.Range("J126").Value = 2770
$ws.Range("K126").Value = 4440
$ws.Range("L126").Value = 8310
$ws.Range("M126").Value = -1970
$ws.Range("N126").Value = -13250

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4073.7693
$ws.Range("I62").Value = 3619.875
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 3619.875
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -2995.875
$ws.Range("N62").Value = -6048
# Row 65
$ws.Range("H65").Value = 4073.7693
$ws.Range("I65").Value = 3619.875
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 18099.375
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -14979.375
$ws.Range("N65").Value = -30240
# Row 126
$ws.Range("H126").Value = 385700.47
$ws.Range("I126").Value = 500755.1
$ws.Range("J126").Value = 2185
$ws.Range("K126").Value = 1502265.3
$ws.Range("L126").Value = 6555
$ws.Range("M126").Value = -1499795.3
$ws.Range("N126").Value = -11495
